$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing the existing rows 7 and 8 down to 8 and 9.
$ws.Range("A7:F7").Insert()

# Copy the style of the row above (row 6) into the newly inserted row 7 so the
# new cells keep the same formatting (centered alignment, etc.) as the rest of the table.
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row 7 with the "All Features" XGBoost entry.
$ws.Range("A7").Value = "XGBoost"
$ws.Range("B7").Value = "All Features"
$ws.Range("C7").Value = "colsample_bytree: 1, learning_rate: 0.3, max_depth: 3, n_estimators: 200, subsample: 0.8"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.9093798223702472
$ws.Range("F7").Value = 0.9012760241773002
